{"js": "// Replace the three-digit x one-digit multiplication expressions in the\n// document's tables with a new set of expressions, per the commit diff.\nconst replacements = [\n  [\"599\u00d74=2396\", \"601\u00d72=1202\"],\n  [\"883\u00d77=6181\", \"381\u00d78=3048\"],\n  [\"269\u00d77=1883\", \"464\u00d73=1392\"],\n  [\"112\u00d73=336\", \"447\u00d78=3576\"],\n  [\"435\u00d78=3480\", \"296\u00d77=2072\"],\n  [\"456\u00d79=4104\", \"825\u00d74=3300\"],\n  [\"857\u00d74=3428\", \"407\u00d73=1221\"],\n  [\"810\u00d74=3240\", \"452\u00d72=904\"],\n  [\"890\u00d78=7120\", \"294\u00d75=1470\"],\n  [\"625\u00d74=2500\", \"199\u00d76=1194\"],\n  [\"474\u00d75=2370\", \"589\u00d75=2945\"],\n  [\"560\u00d75=2800\", \"122\u00d75=610\"],\n  [\"245\u00d78=1960\", \"162\u00d78=1296\"],\n  [\"721\u00d73=2163\", \"299\u00d76=1794\"],\n  [\"210\u00d79=1890\", \"781\u00d76=4686\"],\n  [\"144\u00d79=1296\", \"962\u00d75=4810\"],\n  [\"890\u00d79=8010\", \"952\u00d73=2856\"],\n  [\"621\u00d76=3726\", \"795\u00d72=1590\"],\n  [\"196\u00d73=588\", \"970\u00d77=6790\"],\n  [\"148\u00d76=888\", \"974\u00d78=7792\"],\n  [\"190\u00d76=1140\", \"669\u00d76=4014\"],\n  [\"865\u00d73=2595\", \"588\u00d75=2940\"],\n  [\"222\u00d75=1110\", \"595\u00d75=2975\"],\n  [\"248\u00d73=744\", \"109\u00d77=763\"],\n  [\"737\u00d76=4422\", \"864\u00d77=6048\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication expressions in the\n# document's tables with a new set of expressions, per the commit diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"599\u00d74=2396\", \"601\u00d72=1202\"),\n    @(\"883\u00d77=6181\", \"381\u00d78=3048\"),\n    @(\"269\u00d77=1883\", \"464\u00d73=1392\"),\n    @(\"112\u00d73=336\", \"447\u00d78=3576\"),\n    @(\"435\u00d78=3480\", \"296\u00d77=2072\"),\n    @(\"456\u00d79=4104\", \"825\u00d74=3300\"),\n    @(\"857\u00d74=3428\", \"407\u00d73=1221\"),\n    @(\"810\u00d74=3240\", \"452\u00d72=904\"),\n    @(\"890\u00d78=7120\", \"294\u00d75=1470\"),\n    @(\"625\u00d74=2500\", \"199\u00d76=1194\"),\n    @(\"474\u00d75=2370\", \"589\u00d75=2945\"),\n    @(\"560\u00d75=2800\", \"122\u00d75=610\"),\n    @(\"245\u00d78=1960\", \"162\u00d78=1296\"),\n    @(\"721\u00d73=2163\", \"299\u00d76=1794\"),\n    @(\"210\u00d79=1890\", \"781\u00d76=4686\"),\n    @(\"144\u00d79=1296\", \"962\u00d75=4810\"),\n    @(\"890\u00d79=8010\", \"952\u00d73=2856\"),\n    @(\"621\u00d76=3726\", \"795\u00d72=1590\"),\n    @(\"196\u00d73=588\", \"970\u00d77=6790\"),\n    @(\"148\u00d76=888\", \"974\u00d78=7792\"),\n    @(\"190\u00d76=1140\", \"669\u00d76=4014\"),\n    @(\"865\u00d73=2595\", \"588\u00d75=2940\"),\n    @(\"222\u00d75=1110\", \"595\u00d75=2975\"),\n    @(\"248\u00d73=744\", \"109\u00d77=763\"),\n    @(\"737\u00d76=4422\", \"864\u00d77=6048\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
